$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep text formatting so numeric-looking strings
# (e.g. "575.65", "10.10", "0.0000247") are not coerced into numbers
# and so leading/trailing zeros and thousand-dot formatting survive.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.159.25'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.043.91'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '575.65'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -1.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '169.17'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +3.97%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.041.59'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.521'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.153'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.480'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +5.35%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000247'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -3.18%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.98'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +6.78%  '
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.229.56'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.554.98'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.19%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.23'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +4.12%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.047.24'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.23%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.34'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +17.18%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '467.23'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +2.65%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.706'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +2.58%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.48'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '83.33'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.97'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +5.52%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.10'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -4.14%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.40'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +3.55%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +3.00%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.63'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +7.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0₃0998'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -5.78%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '28.25'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.62%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.84'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '48.70'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +11.02%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.04'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -6.86%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '49.50'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.309'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.121'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.85'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -5.72%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.63'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +2.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0359'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.54%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '384.37'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -3.32%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.744.51'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -2.00%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '134.28'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.92'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +4.18%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.23'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +3.47%  '
